$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Good Morning" (the string previously shown in E8) is replaced with "GIT UPDATE".
# This also drops the now-unused "Good Morning" shared string and appends the
# new "GIT UPDATE" string at the end of the shared-strings table, which shifts
# every other string's index down by one (matching the recorded diff).
$ws.Range("E8").Value = "GIT UPDATE"

# Move / record the active selection on the sheet to E8.
$ws.Range("E8").Select()
